# Applies the "model modifications, updated meta" change to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тесты")

# Turn on wrap text for column H (results); width stays the same.
$ws.Range("H1:H31").WrapText = $true

# Row heights for rows 3 and 4 (new entry plus neighboring existing one) to fit wrapped text.
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 60

# Fill in the new experiment row (row 4) with the Resnet18 results.
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 25
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 6).Value = "Эксперименты с Resnet18"
$ws.Cells.Item(4, 7).Value = "параметры теста 1"
$ws.Cells.Item(4, 8).Value = "Train IoU: 0.4, Val IoU: 0.32. Первые относительно рабочие результаты. Точность всё ещё крайне низкая, но можно дорабатывать отсюда."

# Update selection to reflect the author's last active cell.
$ws.Range("I7").Select()

$wb.Save()
